$d = $word.ActiveDocument

# --- 1. documento.referencia paragraph: drop gramStart/gramEnd proofErr markers ---
$rng = $d.Content
$rng.Find.Execute("documento.referencia") | Out-Null
$p = $rng.Paragraphs(1)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00CC4399" w:rsidRPr="00F92B49" w:rsidRDefault="00CC4399" w:rsidP="00186039"><w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Arial Narrow" w:eastAsia="Calibri" w:hAnsi="Arial Narrow"/><w:b/><w:bCs/><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial Narrow" w:eastAsia="Calibri" w:hAnsi="Arial Narrow"/><w:b/><w:bCs/><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">+++= </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial Narrow" w:eastAsia="Calibri" w:hAnsi="Arial Narrow"/><w:b/><w:bCs/><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>documento.referencia</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial Narrow" w:eastAsia="Calibri" w:hAnsi="Arial Narrow"/><w:b/><w:bCs/><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>+++</w:t></w:r></w:p>')

# --- 2. documento.factura paragraph: drop gramStart/gramEnd proofErr markers ---
$rng = $d.Content
$rng.Find.Execute("documento.factura") | Out-Null
$p = $rng.Paragraphs(1)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00CC4399" w:rsidRPr="00F92B49" w:rsidRDefault="00CC4399" w:rsidP="00186039"><w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Arial Narrow" w:eastAsia="Calibri" w:hAnsi="Arial Narrow"/><w:b/><w:bCs/><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial Narrow" w:eastAsia="Calibri" w:hAnsi="Arial Narrow"/><w:b/><w:bCs/><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">+++= </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial Narrow" w:eastAsia="Calibri" w:hAnsi="Arial Narrow"/><w:b/><w:bCs/><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>documento.factura</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial Narrow" w:eastAsia="Calibri" w:hAnsi="Arial Narrow"/><w:b/><w:bCs/><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>+++</w:t></w:r></w:p>')

# --- 3. documento.cuantia paragraph: drop gramStart/gramEnd + the _GoBack bookmark ---
$rng = $d.Content
$rng.Find.Execute("documento.cuantia") | Out-Null
$p = $rng.Paragraphs(1)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00CC4399" w:rsidRPr="006F0C80" w:rsidRDefault="00CC4399" w:rsidP="00CC4399"><w:pPr><w:pStyle w:val="FR1"/><w:spacing w:line="480" w:lineRule="auto"/><w:ind w:right="2"/><w:rPr><w:rFonts w:ascii="Arial Narrow" w:hAnsi="Arial Narrow" w:cs="Arial"/><w:b/><w:color w:val="000000"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr></w:pPr><w:r w:rsidRPr="006F0C80"><w:rPr><w:rFonts w:ascii="Arial Narrow" w:hAnsi="Arial Narrow" w:cs="Arial"/><w:b/><w:color w:val="000000"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr><w:t xml:space="preserve">CUANTÍA DE LA VENTA: </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial Narrow" w:hAnsi="Arial Narrow" w:cs="Arial"/><w:b/><w:color w:val="000000"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr><w:t xml:space="preserve">+++= </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial Narrow" w:hAnsi="Arial Narrow" w:cs="Arial"/><w:b/><w:color w:val="000000"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr><w:t>documento.cuantia</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial Narrow" w:hAnsi="Arial Narrow" w:cs="Arial"/><w:b/><w:color w:val="000000"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr><w:t>+++</w:t></w:r></w:p>')

# --- 4. documento.copias paragraph: drop gramStart/gramEnd, retarget tail run text,
#        append a new run holding the trailing en-dash ---
$rng = $d.Content
$rng.Find.Execute("documento.copias") | Out-Null
$p = $rng.Paragraphs(1)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00CC4399" w:rsidRPr="006F0C80" w:rsidRDefault="00CC4399" w:rsidP="00CC4399"><w:pPr><w:pStyle w:val="FR1"/><w:spacing w:line="480" w:lineRule="auto"/><w:ind w:right="2"/><w:rPr><w:rFonts w:ascii="Arial Narrow" w:hAnsi="Arial Narrow" w:cs="Arial"/><w:b/><w:color w:val="000000"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr></w:pPr><w:r w:rsidRPr="006F0C80"><w:rPr><w:rFonts w:ascii="Arial Narrow" w:hAnsi="Arial Narrow" w:cs="Arial"/><w:b/><w:color w:val="000000"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr><w:t xml:space="preserve">DI: </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial Narrow" w:hAnsi="Arial Narrow" w:cs="Arial"/><w:b/><w:color w:val="000000"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr><w:t xml:space="preserve">+++= </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial Narrow" w:hAnsi="Arial Narrow" w:cs="Arial"/><w:b/><w:color w:val="000000"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr><w:t>documento.copias</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial Narrow" w:hAnsi="Arial Narrow" w:cs="Arial"/><w:b/><w:color w:val="000000"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr><w:t>+++</w:t></w:r><w:r w:rsidRPr="006F0C80"><w:rPr><w:rFonts w:ascii="Arial Narrow" w:hAnsi="Arial Narrow" w:cs="Arial"/><w:b/><w:color w:val="000000"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr><w:t xml:space="preserve"> COPIAS. </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial Narrow" w:hAnsi="Arial Narrow" w:cs="Arial"/><w:b/><w:color w:val="000000"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr><w:t>–</w:t></w:r></w:p>')

# --- 5. insert the brand-new notary paragraph (carrying the relocated _GoBack bookmark)
#        right before the trailing empty paragraph, leaving that empty paragraph intact ---
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="FR1"/><w:spacing w:line="480" w:lineRule="auto"/><w:ind w:right="2"/><w:rPr><w:rFonts w:ascii="Arial Narrow" w:hAnsi="Arial Narrow" w:cs="Arial"/><w:b/><w:color w:val="000000"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial Narrow" w:hAnsi="Arial Narrow"/><w:color w:val="000000"/><w:sz w:val="25"/><w:szCs w:val="25"/></w:rPr><w:t xml:space="preserve">En la Ciudad de Portoviejo, Capital de la Provincia de Manabí, República del Ecuador, hoy </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial Narrow" w:hAnsi="Arial Narrow"/><w:color w:val="000000"/><w:sz w:val="25"/><w:szCs w:val="25"/></w:rPr><w:t xml:space="preserve">+++= </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial Narrow" w:hAnsi="Arial Narrow"/><w:color w:val="000000"/><w:sz w:val="25"/><w:szCs w:val="25"/></w:rPr><w:t>documento.fecha</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial Narrow" w:hAnsi="Arial Narrow"/><w:color w:val="000000"/><w:sz w:val="25"/><w:szCs w:val="25"/></w:rPr><w:t>+++</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial Narrow" w:hAnsi="Arial Narrow"/><w:color w:val="000000"/><w:sz w:val="25"/><w:szCs w:val="25"/></w:rPr><w:t xml:space="preserve">, ante mi </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial Narrow" w:hAnsi="Arial Narrow"/><w:color w:val="000000"/><w:sz w:val="25"/><w:szCs w:val="25"/></w:rPr><w:t xml:space="preserve">+++= </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial Narrow" w:hAnsi="Arial Narrow"/><w:color w:val="000000"/><w:sz w:val="25"/><w:szCs w:val="25"/></w:rPr><w:t>documento.notario_leyenda</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial Narrow" w:hAnsi="Arial Narrow"/><w:color w:val="000000"/><w:sz w:val="25"/><w:szCs w:val="25"/></w:rPr><w:t>+++</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial Narrow" w:hAnsi="Arial Narrow"/><w:color w:val="000000"/><w:sz w:val="25"/><w:szCs w:val="25"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial Narrow" w:hAnsi="Arial Narrow"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="25"/><w:szCs w:val="25"/><w:lang w:val="es-EC"/></w:rPr><w:t>comparece</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial Narrow" w:hAnsi="Arial Narrow"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="25"/><w:szCs w:val="25"/><w:lang w:val="es-EC"/></w:rPr><w:t>n: por una parte</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial Narrow" w:hAnsi="Arial Narrow"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="25"/><w:szCs w:val="25"/><w:lang w:val="es-EC"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>')
